# =====================================================================
# Edit script: "add get fast result with tables"
#
# Summary of changes (derived from the OOXML diff):
#  1. Insert a new worksheet "Аканское" right after "Архангельское"
#     (2nd tab), with a small 2-column summary table.
#  2. On "Матросовское" delete the (empty) leading column A so the
#     data shifts one column to the left.
#  3. Remove the now-vestigial, empty, styled cells in column T
#     (rows 2-10) on "Архангельское" (T1 is kept).
#  4. Misc. view/selection/zoom housekeeping on several sheets,
#     including moving the active tab to "Байданкинское".
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New sheet "Аканское" inserted after "Архангельское"
# ---------------------------------------------------------------------
$sheetArk = $wb.Worksheets.Item("Архангельское")
$akan = $wb.Worksheets.Add($null, $sheetArk)
$akan.Name = "Аканское"

$akan.Range("A1").Value = "Месторождение"
$akan.Range("B1").Value = "Местоположение"
$akan.Range("A2").Value = "Аканское месторождение"
$akan.Range("B2").Value = "нурлатский район республика татарстан"

$akanHeader = $akan.Range("A1:B1")
$akanHeader.Font.Bold = $true
$akanHeader.HorizontalAlignment = -4108
$akanHeader.VerticalAlignment = -4160
$akanHeader.Borders.LineStyle = 1

$akan.Columns.Item(1).ColumnWidth = 25.44140625
$akan.Columns.Item(2).ColumnWidth = 38.88671875

# ---------------------------------------------------------------------
# 2. "Матросовское": drop the empty first column, shifting left
# ---------------------------------------------------------------------
$matr = $wb.Worksheets.Item("Матросовское")
$matr.Columns.Item(1).Delete()

$matr.Activate()
$matr.Range("D8").Select()

# ---------------------------------------------------------------------
# 3. "Архангельское": clear the leftover styled cells T2:T10
#    (T1 itself is left in place)
# ---------------------------------------------------------------------
$ark = $wb.Worksheets.Item("Архангельское")
$ark.Range("T2:T10").Clear()

$ark.Activate()
$excel.ActiveWindow.Zoom = 69
$ark.Range("A1:R12").Select()
$ark.Range("C12").Activate()

# ---------------------------------------------------------------------
# 4. Misc. selection / zoom bookkeeping on the remaining sheets
# ---------------------------------------------------------------------
$ivin = $wb.Worksheets.Item("Ивинское")
$ivin.Activate()
$ivin.Range("E12").Select()

$gran = $wb.Worksheets.Item("Граничное")
$gran.Activate()
$gran.Cells.Select()
$gran.Range("D8").Activate()

$baid = $wb.Worksheets.Item("Байданкинское")
$baid.Activate()
$baid.Cells.Select()
$baid.Range("C19").Activate()

$sherb = $wb.Worksheets.Item("Щербенское")
$sherb.Activate()
$sherb.Range("C2").Select()

# Final active sheet/tab -> "Байданкинское" (matches workbook activeTab)
$baid.Activate()
